$wb = $excel.ActiveWorkbook

# Sheet "展览" (index 1 / rId1) - source rows
$wsExhibit = $wb.Worksheets.Item(1)
$wsExhibit.Range("F5").Value = 303
$wsExhibit.Range("F7").Value = 1050
$wsExhibit.Range("F9").Value = 545
$wsExhibit.Range("F10").Value = 563
$wsExhibit.Range("F11").Value = 169
$wsExhibit.Range("F12").Value = 13398
$wsExhibit.Range("F14").Value = 17
$wsExhibit.Range("F16").Value = 5519

# Sheet "全部类型" (index 4 / rId4) - combined/all-types view, same records duplicated
$wsAll = $wb.Worksheets.Item(4)
$wsAll.Range("F21").Value = 303
$wsAll.Range("F29").Value = 1050
$wsAll.Range("F31").Value = 545
$wsAll.Range("F32").Value = 563
$wsAll.Range("F33").Value = 169
$wsAll.Range("F34").Value = 13398
$wsAll.Range("F36").Value = 17
$wsAll.Range("F39").Value = 5519
